$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/local-number-of-units-per-service"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
# "Fixed Value" for the Extension.url row (row 5) mirrors the same URL change
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/local-number-of-units-per-service"
# "Constraint(s)" for the root Extension row (row 2) is cleared
$elements.Range("AI2").Value = ""
